$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff": rows 4-7 on the zh-cn and de-de sheets
# move from "low" priority to "ht" (handed-off) priority, and their
# "Latest Handoff Datetime" timestamps are refreshed to the handoff time.
# The Overview sheet's "Latest HO Xliff Generate Date" column for the
# same files tracks the de-de handoff timestamp, so it is refreshed too.

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

for ($r = 4; $r -le 7; $r++) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-16 16:27:52"

    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-16 16:27:57"

    $overview.Cells.Item($r, 7).Value = "2016-08-16 16:27:57"
}
